# Update BoM row for U105 U106 (ACS712 current sensor) from the ±5A/05B
# variant to the ±20A/20A variant, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

$ws.Range("C29").Value = "ACS712xLCTR-20A"
$ws.Range("D29").Value = "ACS712ELCTR-20A"
$ws.Range("E29").Value = "ACS712ELCTR-20A-T"
$ws.Range("I29").Value = "https://jlcpcb.com/partdetail/11225-ACS712ELCTR_20AT/C10681"
$ws.Range("K29").Value = "±20A Bidirectional Hall-Effect Current Sensor, +5.0V supply, 100mV/A, SOIC-8"
